# Generate Report for Handoff
# Replaces the stale e2e test-file identifiers/hashes with the freshly
# generated ones, refreshes the associated timestamps, and clears out the
# (not-yet-existent) handback file/date columns for the zh-cn and de-de
# localization sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "98e559ef-0f65-455a-b0c2-06267c490dcb"
$newGuid = "b1556eff-bd22-4bed-aa91-97f977afd908"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-01 03:10:58"

$overviewUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1754fd32a6397d4313a96f96e8f972a375af70e3/e2e/$oldGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewUrl, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.cf5aec950a6579089e769db15f2fa67bdfea96b9.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 03:10:53"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I2").Style = "Normal"

$zhCnUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1754fd32a6397d4313a96f96e8f972a375af70e3/e2e/$oldGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnUrl, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.cf5aec950a6579089e769db15f2fa67bdfea96b9.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 03:10:58"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I2").Style = "Normal"

$deDeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1754fd32a6397d4313a96f96e8f972a375af70e3/e2e/$oldGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeUrl, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
